$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (style) from H1 into the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-17
$data = @(
    @(5,6),
    @(9,9),
    @(4,6),
    @(7,8),
    @(9,9),
    @(8,9),
    @(9,9),
    @(7,7),
    @(9,9),
    @(8,9),
    @(7,8),
    @(8,9),
    @(9,9),
    @(7,8),
    @(4,6),
    @(4,4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
